$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALERT")

# New row (row 6) for ALERT_005, following the same pattern as the
# existing ALERT_00x rows (2-5).
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "ALERT_005"

$ws.Range("C6").Value = $ws.Range("C5").Value()
$ws.Range("C6").NumberFormat = $ws.Range("C5").NumberFormat()
$ws.Range("C6").Font.Name = $ws.Range("C5").Font.Name()
$ws.Range("C6").Font.Size = $ws.Range("C5").Font.Size()

$ws.Range("D6").Value = $ws.Range("D5").Value()

$ws.Range("E6").Value = "Customer_5"

$ws.Range("F6").Value = "Customer Name|CIF Number|Alert Heading|Alert Content|User Name|Date Added / Amended"
$ws.Range("F6").NumberFormat = $ws.Range("F5").NumberFormat()
$ws.Range("F6").Font.Name = $ws.Range("F5").Font.Name()
$ws.Range("F6").Font.Size = $ws.Range("F5").Font.Size()

# Keep the active-cell selection pointing one row below the new data,
# mirroring the prior "next empty row" selection behaviour.
$ws.Range("F7").Select()
